# Appends the 7 newly-scraped OLX listings (rows 47-53) to the
# PODSUMOWANIE sheet running log, mirroring the style of the
# existing rows for the same profile/listing further up the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PODSUMOWANIE")

# Row 47: same listing as row 7, re-checked at the newer timestamp.
$ws.Range("A7:H7").Copy($ws.Range("A47:H47"))
$ws.Range("A47").Value = "2026-02-17 16:33:37"
$ws.Range("B47").Value = "poqui"
$ws.Range("C47").Value = "Świeżo wykończone mieszkanie z dużym balkonem - Ponikwoda"
$ws.Range("D47").Value = 2299
# Force literal text for the date so Excel does not reinterpret
# an unambiguous dd.mm.yyyy string as a date serial, then restore
# the normal (General) look of the column from the template cell.
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "19.01.2026"
$ws.Range("E7").Copy() | Out-Null
$ws.Range("E47").PasteSpecial(-4122) | Out-Null
$ws.Range("F47").Value = 29
$ws.Range("G47").Value = "https://www.olx.pl/d/oferta/swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR.html"
$ws.Range("H47").Value = "swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR"

# Row 48: same listing as row 8, re-checked at the newer timestamp.
$ws.Range("A8:H8").Copy($ws.Range("A48:H48"))
$ws.Range("A48").Value = "2026-02-17 16:33:37"
$ws.Range("B48").Value = "poqui"
$ws.Range("C48").Value = "Kawalerka po remoncie z funkcjonalną antresolą - ul. Jana Sawy"
$ws.Range("D48").Value = 2499
# Force literal text for the date so Excel does not reinterpret
# an unambiguous dd.mm.yyyy string as a date serial, then restore
# the normal (General) look of the column from the template cell.
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "28.10.2025"
$ws.Range("E8").Copy() | Out-Null
$ws.Range("E48").PasteSpecial(-4122) | Out-Null
$ws.Range("F48").Value = 112
$ws.Range("G48").Value = "https://www.olx.pl/d/oferta/kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger.html"
$ws.Range("H48").Value = "kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger"

# Row 49: same listing as row 9, re-checked at the newer timestamp.
$ws.Range("A9:H9").Copy($ws.Range("A49:H49"))
$ws.Range("A49").Value = "2026-02-17 16:33:37"
$ws.Range("B49").Value = "poqui"
$ws.Range("C49").Value = "Przytulny pokój blisko Politechniki – ul. Przytulna"
$ws.Range("D49").Value = 599
# Force literal text for the date so Excel does not reinterpret
# an unambiguous dd.mm.yyyy string as a date serial, then restore
# the normal (General) look of the column from the template cell.
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "10.10.2025"
$ws.Range("E9").Copy() | Out-Null
$ws.Range("E49").PasteSpecial(-4122) | Out-Null
$ws.Range("F49").Value = 130
$ws.Range("G49").Value = "https://www.olx.pl/d/oferta/przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz.html"
$ws.Range("H49").Value = "przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz"

# Row 50: same listing as row 11, re-checked at the newer timestamp.
$ws.Range("A11:H11").Copy($ws.Range("A50:H50"))
$ws.Range("A50").Value = "2026-02-17 16:33:37"
$ws.Range("B50").Value = "pokojewlublinie"
$ws.Range("C50").Value = "WOLNY OD ZARAZ! Pokój jedynka, ul. Romanowskiego 58"
$ws.Range("D50").Value = 58640
# Force literal text for the date so Excel does not reinterpret
# an unambiguous dd.mm.yyyy string as a date serial, then restore
# the normal (General) look of the column from the template cell.
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "11.08.2025"
$ws.Range("E11").Copy() | Out-Null
$ws.Range("E50").PasteSpecial(-4122) | Out-Null
$ws.Range("F50").Value = 190
$ws.Range("G50").Value = "https://www.olx.pl/d/oferta/wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm.html"
$ws.Range("H50").Value = "wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm"

# Row 51: same listing as row 12, re-checked at the newer timestamp.
$ws.Range("A12:H12").Copy($ws.Range("A51:H51"))
$ws.Range("A51").Value = "2026-02-17 16:33:37"
$ws.Range("B51").Value = "pokojewlublinie"
$ws.Range("C51").Value = "WOLNY OD ZARAZ! Super lokalizacja, blisko centrum, ul. Paganiniego 12"
$ws.Range("D51").Value = 12640
# Force literal text for the date so Excel does not reinterpret
# an unambiguous dd.mm.yyyy string as a date serial, then restore
# the normal (General) look of the column from the template cell.
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "19.01.2026"
$ws.Range("E12").Copy() | Out-Null
$ws.Range("E51").PasteSpecial(-4122) | Out-Null
$ws.Range("F51").Value = 29
$ws.Range("G51").Value = "https://www.olx.pl/d/oferta/wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc.html"
$ws.Range("H51").Value = "wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc"

# Row 52: same listing as row 13, re-checked at the newer timestamp.
$ws.Range("A13:H13").Copy($ws.Range("A52:H52"))
$ws.Range("A52").Value = "2026-02-17 16:33:37"
$ws.Range("B52").Value = "dawnypatron"
$ws.Range("C52").Value = "Ładny pokój jednoosobowy. Wynajmę duży pokój w centrum. ul Niecała 4."
$ws.Range("D52").Value = 730
# Force literal text for the date so Excel does not reinterpret
# an unambiguous dd.mm.yyyy string as a date serial, then restore
# the normal (General) look of the column from the template cell.
$ws.Range("E52").NumberFormat = "@"
$ws.Range("E52").Value = "20.09.2024"
$ws.Range("E13").Copy() | Out-Null
$ws.Range("E52").PasteSpecial(-4122) | Out-Null
$ws.Range("F52").Value = 515
$ws.Range("G52").Value = "https://www.olx.pl/d/oferta/ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM.html"
$ws.Range("H52").Value = "ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM"

# Row 53: same listing as row 14, re-checked at the newer timestamp.
$ws.Range("A14:H14").Copy($ws.Range("A53:H53"))
$ws.Range("A53").Value = "2026-02-17 16:33:37"
$ws.Range("B53").Value = "dawnypatron"
$ws.Range("C53").Value = "Mam do wynajęcia pokój dla os. pracującej lub studenta. Narutowicza 14"
$ws.Range("D53").Value = 14690
# Force literal text for the date so Excel does not reinterpret
# an unambiguous dd.mm.yyyy string as a date serial, then restore
# the normal (General) look of the column from the template cell.
$ws.Range("E53").NumberFormat = "@"
$ws.Range("E53").Value = "05.12.2025"
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E53").PasteSpecial(-4122) | Out-Null
$ws.Range("F53").Value = 74
$ws.Range("G53").Value = "https://www.olx.pl/d/oferta/mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv.html"
$ws.Range("H53").Value = "mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv"

$excel.CutCopyMode = $false
Write-Output "OLX monitor: appended rows 47-53 (2026-02-17 16:33:37 check)."
